$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-derived values for rows 2-7 (columns D through T).
# Column order: D,E,F,G,H,I,J,K,L,M,N,O,P,Q,R,S,T

$data = @{
    2 = @("ECs",    3, 1,                  7.514794999999999,  22.544385, 0.977669497583861,  0.977669497583861,  3, 1, 0.491263,            1.473789,  0.08507192955174298, 0.08507192955174299, 3.691740736085,    33.225666624765,    0.08317223062334217, 0.08317223062334218)
    3 = @("FAPs",   3, 1,                  7.514794999999999,  22.544385, 0.977669497583861,  0.977669497583861,  3, 1, 4.669265333333333,  14.007796, 0.8085758778815603,  0.8085758778815605,  35.08857178060666, 315.7971460254599,  0.7905199722868944,  0.7905199722868946)
    4 = @("MuSCs",  3, 1,                  7.514794999999999,  22.544385, 0.977669497583861,  0.977669497583861,  3, 1, 0.6141496666666667, 1.842449,  0.1063521925666967,  0.1063521925666967,  4.615208844318333, 41.536879598865,    0.1039772946736244,  0.1039772946736244)
    5 = @("ECs",    2, 0.6666666666666666, 0.171642,            0.514926,  0.02233050241613897, 0.02233050241613898, 3, 1, 0.491263,         1.473789,  0.08507192955174298, 0.08507192955174299, 0.08432136384599999, 0.758892274614,   0.001899698928400801, 0.001899698928400802)
    6 = @("FAPs",   2, 0.6666666666666666, 0.171642,            0.514926,  0.02233050241613897, 0.02233050241613898, 3, 1, 4.669265333333333, 14.007796, 0.8085758778815603,  0.8085758778815605,  0.8014420403439998, 7.212978363095999, 0.01805590559466587, 0.01805590559466588)
    7 = @("MuSCs",  2, 0.6666666666666666, 0.171642,            0.514926,  0.02233050241613897, 0.02233050241613898, 3, 1, 0.6141496666666667, 1.842449, 0.1063521925666967,  0.1063521925666967,  0.105413877086,    0.948724893774,     0.002374897893072297, 0.002374897893072297)
}

foreach ($r in 2..7) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 4 + $i   # column D is index 4
        $ws.Cells.Item($r, $col).Value = $vals[$i]
    }
}
